$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.824.83"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").Value = "1.887.70"
$ws.Range("E3").Value = "  -0.40%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7529"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.31"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.67%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9998"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3121"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.31"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.79%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07125"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08486"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.91%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7598"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.72%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.875.66"
$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.358"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.58%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.38"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.147"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.36%  "

$ws.Range("D17").Value = "29.854.91"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("E18").Value = "  -1.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "243.66"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007802"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "2.149.45"
$ws.Range("E21").Value = "  +1.02%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9991"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.014"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9998"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1596"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.11%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.366"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "162.27"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.73"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.028"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.494"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.539"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.511"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.51%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.129"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"

$ws.Range("E34").Value = "  -2.75%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7501"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.003"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.707"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.86%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01950"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.772"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.90%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4459"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.37%  "

$ws.Range("D43").Value = "1.090.07"
$ws.Range("E43").Value = "  -3.70%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.57"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.33%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8576"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9997"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.724"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.59%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "102.45"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.859"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.052"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.04%  "

$ws.Range("D51").Value = "2.031.73"
$ws.Range("E51").Value = "  -0.39%  "
